$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1012.6667
$ws.Range("I39").Value = 301.25
$ws.Range("J39").Value = 1581.8
$ws.Range("K39").Value = 903.75
$ws.Range("L39").Value = 4745.4
$ws.Range("M39").Value = -607.75
$ws.Range("N39").Value = -5337.4

$ws.Range("H40").Value = 1050.6571
$ws.Range("I40").Value = 1080.7931
$ws.Range("J40").Value = 905
$ws.Range("K40").Value = 1080.7931
$ws.Range("L40").Value = 905
$ws.Range("M40").Value = -905.7931000000001
$ws.Range("N40").Value = -1255

$ws.Range("H64").Value = 3550
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3550
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3550
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4046

$ws.Range("H67").Value = 3550
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3550
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3550
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5266

$ws.Range("H137").Value = 1031.3572
$ws.Range("I137").Value = 788.36365
$ws.Range("J137").Value = 1922.3334
$ws.Range("K137").Value = 2365.09095
$ws.Range("L137").Value = 5767.0002
$ws.Range("M137").Value = 184.9090500000002
$ws.Range("N137").Value = -10867.0002

$ws.Range("H138").Value = 2082.262
$ws.Range("I138").Value = 2125.2122
$ws.Range("K138").Value = 6375.6366
$ws.Range("M138").Value = -1235.6366

$ws.Range("H139").Value = 50255.445
$ws.Range("J139").Value = 50255.445
$ws.Range("L139").Value = 50255.445
$ws.Range("N139").Value = -60535.445

$ws.Range("H140").Value = 56884.668
$ws.Range("J140").Value = 56884.668
$ws.Range("L140").Value = 56884.668
$ws.Range("N140").Value = -67244.66800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 794356.4399999999
$ws.Range("I2").Value = 2778140
$ws.Range("J2").Value = 843
$ws.Range("K2").Value = 2778140
$ws.Range("L2").Value = 843
$ws.Range("M2").Value = -2778027
$ws.Range("N2").Value = -1069

$ws.Range("H116").Value = 794356.4399999999
$ws.Range("I116").Value = 2778140
$ws.Range("J116").Value = 843
$ws.Range("K116").Value = 2778140
$ws.Range("L116").Value = 843
$ws.Range("M116").Value = -2775846
$ws.Range("N116").Value = -5431

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 794356.4399999999
$ws.Range("I3").Value = 2778140
$ws.Range("J3").Value = 843
$ws.Range("K3").Value = 2778140
$ws.Range("L3").Value = 843
$ws.Range("M3").Value = -2778026
$ws.Range("N3").Value = -1071

$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 10000
$ws.Range("K35").Value = 10000
$ws.Range("M35").Value = -9690

$ws.Range("H134").Value = 8952.84
$ws.Range("I134").Value = 10620.1875
$ws.Range("J134").Value = 5988.6665
$ws.Range("K134").Value = 31860.5625
$ws.Range("L134").Value = 17965.9995
$ws.Range("M134").Value = -29325.5625
$ws.Range("N134").Value = -23035.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H31").Value = 1778.4762
$ws.Range("I31").Value = 1233.4615
$ws.Range("K31").Value = 1233.4615
$ws.Range("M31").Value = -938.4614999999999

$ws.Range("H34").Value = 1778.4762
$ws.Range("I34").Value = 1233.4615
$ws.Range("K34").Value = 1233.4615
$ws.Range("M34").Value = -1031.4615

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 1397.9584
$ws.Range("I132").Value = 784.439
$ws.Range("K132").Value = 2353.317
$ws.Range("M132").Value = 176.683

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 277.2
$ws.Range("I6").Value = 296.5
$ws.Range("K6").Value = 889.5
$ws.Range("M6").Value = -776.5

$ws.Range("H8").Value = 43.25
$ws.Range("I8").Value = 43.25
$ws.Range("K8").Value = 129.75
$ws.Range("M8").Value = 9.25

$ws.Range("H11").Value = 752.5
$ws.Range("I11").Value = 752.5
$ws.Range("K11").Value = 2257.5
$ws.Range("M11").Value = -2117.5

$ws.Range("H19").Value = 2000
$ws.Range("I19").Value = 2000
$ws.Range("K19").Value = 6000
$ws.Range("M19").Value = -5826

$ws.Range("H22").Value = 2000.5
$ws.Range("I22").Value = 2000.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 6001.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5832.5
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 2000.5
$ws.Range("I27").Value = 2000.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 6001.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -5899.5
$ws.Range("N27").ClearContents()

$ws.Range("H109").Value = 4140.364
$ws.Range("I109").Value = 1344.3334
$ws.Range("J109").Value = 5188.875
$ws.Range("K109").Value = 4033.0002
$ws.Range("L109").Value = 15566.625
$ws.Range("M109").Value = -2993.0002
$ws.Range("N109").Value = -17646.625

$ws.Range("H121").Value = 599.25
$ws.Range("I121").Value = 615.25
$ws.Range("J121").Value = 583.25
$ws.Range("K121").Value = 1845.75
$ws.Range("L121").Value = 1749.75
$ws.Range("M121").Value = -535.75
$ws.Range("N121").Value = -4369.75

$ws.Range("H122").Value = 1210.4375
$ws.Range("I122").Value = 1001
$ws.Range("J122").Value = 1258.7693
$ws.Range("K122").Value = 9009
$ws.Range("L122").Value = 11328.9237
$ws.Range("M122").Value = -6559
$ws.Range("N122").Value = -16228.9237

$ws.Range("H131").Value = 15021.036
$ws.Range("I131").Value = 366.63635
$ws.Range("J131").Value = 18603.223
$ws.Range("K131").Value = 1099.90905
$ws.Range("L131").Value = 55809.66900000001
$ws.Range("M131").Value = 3940.09095
$ws.Range("N131").Value = -65889.66900000001

$ws.Range("H134").Value = 2718.353
$ws.Range("I134").Value = 1278.25
$ws.Range("J134").Value = 3998.4443
$ws.Range("K134").Value = 3834.75
$ws.Range("L134").Value = 11995.3329
$ws.Range("M134").Value = 1235.25
$ws.Range("N134").Value = -22135.3329

$ws.Range("H136").Value = 3159.8
$ws.Range("I136").Value = 3159.8
$ws.Range("K136").Value = 9479.400000000001
$ws.Range("M136").Value = -4379.400000000001

$ws.Range("H140").Value = 2218.963
$ws.Range("I140").Value = 1341.1428
$ws.Range("J140").Value = 2526.2
$ws.Range("K140").Value = 4023.4284
$ws.Range("L140").Value = 7578.599999999999
$ws.Range("M140").Value = 1156.5716
$ws.Range("N140").Value = -17938.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 26010.5
$ws.Range("I26").Value = 14000
$ws.Range("K26").Value = 14000
$ws.Range("M26").Value = -13720

$ws.Range("H50").Value = 26010.5
$ws.Range("I50").Value = 14000
$ws.Range("K50").Value = 14000
$ws.Range("M50").Value = -13502

$ws.Range("H102").Value = 2519.5625
$ws.Range("I102").Value = 2522.4285
$ws.Range("K102").Value = 2522.4285
$ws.Range("M102").Value = -900.4285

$ws.Range("H127").Value = 36605.668
$ws.Range("J127").Value = 36605.668
$ws.Range("L127").Value = 36605.668
$ws.Range("N127").Value = -46525.668

$ws.Range("H139").Value = 60162.5
$ws.Range("J139").Value = 60162.5
$ws.Range("L139").Value = 60162.5
$ws.Range("N139").Value = -70442.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 350481.34
$ws.Range("J64").Value = 25722.5
$ws.Range("L64").Value = 25722.5
$ws.Range("N64").Value = -26172.5

$ws.Range("H67").Value = 350481.34
$ws.Range("J67").Value = 25722.5
$ws.Range("L67").Value = 25722.5
$ws.Range("N67").Value = -27282.5

$ws.Range("H122").Value = 5323.8335
$ws.Range("I122").Value = 8269.25
$ws.Range("J122").Value = 3851.125
$ws.Range("K122").Value = 24807.75
$ws.Range("L122").Value = 11553.375
$ws.Range("M122").Value = -22357.75
$ws.Range("N122").Value = -16453.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1256.4
$ws.Range("I113").Value = 760.6667
$ws.Range("K113").Value = 2282.0001
$ws.Range("M113").Value = -112.0001000000002

$ws.Range("H119").Value = 26088.6
$ws.Range("J119").Value = 26088.6
$ws.Range("L119").Value = 26088.6
$ws.Range("N119").Value = -35764.6

$ws.Range("H122").Value = 87505
$ws.Range("I122").Value = 98336.5
$ws.Range("J122").Value = 853
$ws.Range("K122").Value = 295009.5
$ws.Range("L122").Value = 2559
$ws.Range("M122").Value = -292559.5
$ws.Range("N122").Value = -7459

$ws.Range("H136").Value = 9579948
$ws.Range("I136").Value = 13228859
$ws.Range("K136").Value = 39686577
$ws.Range("M136").Value = -39684027

$ws.Range("H139").Value = 69966.664
$ws.Range("L139").Value = 69966.664
$ws.Range("N139").Value = -80246.664
